$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing score for row 18
$ws.Range("C18").Value = 0.715

# Copy the date formatting used by existing rows down onto the new rows
$ws.Range("A18").Copy()
$ws.Range("A19:A21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Add three new rows of results
$ws.Range("A19").Value = 44459
$ws.Range("B19").Value = "model_floodwater_unet_pc_augm_diceloss_5"
$ws.Range("C19").Value = 0.704

$ws.Range("A20").Value = 44459
$ws.Range("B20").Value = "model_floodwater_unet_pc_augm_diceloss 1 +2 + 5"
$ws.Range("C20").Value = 0.718

$ws.Range("A21").Value = 44459
$ws.Range("B21").Value = "model_floodwater_unet_pc_augm_diceloss 1 +2 + 3 + 4 + 5"
$ws.Range("C21").Value = 0.716

# Update the view state to match the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D15").Select()

$wb.Save()
